# Update Data Sources from LFX: refresh the table style used by every
# data-source table in the deck (old GUID -> new GUID).
$oldStyleId = "{C320EC70-FFC7-4695-BD8A-B35883984C6E}"
$newStyleId = "{50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}"

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTable) {
            $tbl = $shp.Table
            $tbl.ApplyStyle($newStyleId)
        }
    }
}
